$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.6594814792829158
$ws.Range("J2").Value = 0.6594814792829158
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.05601
$ws.Range("N2").Value = 0.16803
$ws.Range("O2").Value = 0.02710547761971223
$ws.Range("P2").Value = 0.02710547761971223
$ws.Range("Q2").Value = 0.00331866718
$ws.Range("R2").Value = 0.02986800462
$ws.Range("S2").Value = 0.01787556047731779
$ws.Range("T2").Value = 0.01787556047731779

# Row 3
$ws.Range("I3").Value = 0.6594814792829158
$ws.Range("J3").Value = 0.6594814792829158
$ws.Range("N3").Value = 5.594253
$ws.Range("O3").Value = 0.902427539668559
$ws.Range("P3").Value = 0.9024275396685592
$ws.Range("Q3").Value = 0.1104889830846667
$ws.Range("R3").Value = 0.994400847762
$ws.Range("S3").Value = 0.5951342488062635
$ws.Range("T3").Value = 0.5951342488062635

# Row 4
$ws.Range("I4").Value = 0.6594814792829158
$ws.Range("J4").Value = 0.6594814792829158
$ws.Range("M4").Value = 0.145611
$ws.Range("N4").Value = 0.436833
$ws.Range("O4").Value = 0.07046698271172858
$ws.Range("P4").Value = 0.07046698271172858
$ws.Range("Q4").Value = 0.008627645898000001
$ws.Range("R4").Value = 0.077648813082
$ws.Range("S4").Value = 0.04647166999933441
$ws.Range("T4").Value = 0.04647166999933441

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.030594
$ws.Range("H5").Value = 0.091782
$ws.Range("I5").Value = 0.3405185207170842
$ws.Range("J5").Value = 0.3405185207170842
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.05601
$ws.Range("N5").Value = 0.16803
$ws.Range("O5").Value = 0.02710547761971223
$ws.Range("P5").Value = 0.02710547761971223
$ws.Range("Q5").Value = 0.00171356994
$ws.Range("R5").Value = 0.01542212946
$ws.Range("S5").Value = 0.009229917142394439
$ws.Range("T5").Value = 0.009229917142394441

# Row 6
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.030594
$ws.Range("H6").Value = 0.091782
$ws.Range("I6").Value = 0.3405185207170842
$ws.Range("J6").Value = 0.3405185207170842
$ws.Range("N6").Value = 5.594253
$ws.Range("O6").Value = 0.902427539668559
$ws.Range("P6").Value = 0.9024275396685592
$ws.Range("Q6").Value = 0.057050192094
$ws.Range("R6").Value = 0.513451728846
$ws.Range("S6").Value = 0.3072932908622955
$ws.Range("T6").Value = 0.3072932908622956

# Row 7
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.030594
$ws.Range("H7").Value = 0.091782
$ws.Range("I7").Value = 0.3405185207170842
$ws.Range("J7").Value = 0.3405185207170842
$ws.Range("M7").Value = 0.145611
$ws.Range("N7").Value = 0.436833
$ws.Range("O7").Value = 0.07046698271172858
$ws.Range("P7").Value = 0.07046698271172858
$ws.Range("Q7").Value = 0.004454822934
$ws.Range("R7").Value = 0.040093406406
$ws.Range("S7").Value = 0.02399531271239416
$ws.Range("T7").Value = 0.02399531271239416
